# Auto-applied BRVM data refresh (GitHub Actions update)
# Updates "Recommandations" (A2:G50) and "Top_YTD" (B2:B11) sheets
# with refreshed market figures, rankings and a handful of re-sorted
# rows (the source table is kept sorted by "Variation Totale (%)").

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Sheet "Recommandations": rows 2-50 (header row 1 untouched) ---
# Row 2: SUCRIVOIRE
$ws1.Range("A2").Value = "SUCRIVOIRE"
$ws1.Range("B2").Value = 0
$ws1.Range("C2").Value = 4
$ws1.Range("D2").Value = 3930
$ws1.Range("E2").Value = 985
$ws1.Range("F2").Value = "🟡 Observer"
$ws1.Range("G2").Value = "➖ Neutre"
# Row 3: BRVM - SERVICES PUBLICS
$ws1.Range("A3").Value = "BRVM - SERVICES PUBLICS"
$ws1.Range("B3").Value = 0
$ws1.Range("C3").Value = 8
$ws1.Range("D3").Value = 3364.64
$ws1.Range("E3").Value = 111.7
$ws1.Range("F3").Value = "🟡 Observer"
$ws1.Range("G3").Value = "➖ Neutre"
# Row 4: SAFCA CI
$ws1.Range("A4").Value = "SAFCA CI"
$ws1.Range("B4").Value = 0
$ws1.Range("C4").Value = 4
$ws1.Range("D4").Value = 2730
$ws1.Range("E4").Value = 700
$ws1.Range("F4").Value = "🟡 Observer"
$ws1.Range("G4").Value = "➖ Neutre"
# Row 5: CFAO MOTORS CI
$ws1.Range("A5").Value = "CFAO MOTORS CI"
$ws1.Range("B5").Value = 0
$ws1.Range("C5").Value = 4
$ws1.Range("D5").Value = 2725
$ws1.Range("E5").Value = 675
$ws1.Range("F5").Value = "🟡 Observer"
$ws1.Range("G5").Value = "➖ Neutre"
# Row 6: BRVM - AUTRES SECTEURS
$ws1.Range("A6").Value = "BRVM - AUTRES SECTEURS"
$ws1.Range("B6").Value = 0
$ws1.Range("C6").Value = 4
$ws1.Range("D6").Value = 2651.92
$ws1.Range("E6").Value = 663.72
$ws1.Range("F6").Value = "🟡 Observer"
$ws1.Range("G6").Value = "➖ Neutre"
# Row 7: NEI-CEDA CI
$ws1.Range("A7").Value = "NEI-CEDA CI"
$ws1.Range("B7").Value = 0
$ws1.Range("C7").Value = 4
$ws1.Range("D7").Value = 2355
$ws1.Range("E7").Value = 600
$ws1.Range("F7").Value = "🟡 Observer"
$ws1.Range("G7").Value = "➖ Neutre"
# Row 8: UNIWAX CI
$ws1.Range("A8").Value = "UNIWAX CI"
$ws1.Range("B8").Value = 0
$ws1.Range("C8").Value = 4
$ws1.Range("D8").Value = 2345
$ws1.Range("E8").Value = 595
$ws1.Range("F8").Value = "🟡 Observer"
$ws1.Range("G8").Value = "➖ Neutre"
# Row 9: SETAO CI
$ws1.Range("A9").Value = "SETAO CI"
$ws1.Range("B9").Value = 0
$ws1.Range("C9").Value = 4
$ws1.Range("D9").Value = 2285
$ws1.Range("E9").Value = 575
$ws1.Range("F9").Value = "🟡 Observer"
$ws1.Range("G9").Value = "➖ Neutre"
# Row 10: AIR LIQUIDE CI
$ws1.Range("A10").Value = "AIR LIQUIDE CI"
$ws1.Range("B10").Value = 0
$ws1.Range("C10").Value = 4
$ws1.Range("D10").Value = 2145
$ws1.Range("E10").Value = 535
$ws1.Range("F10").Value = "🟡 Observer"
$ws1.Range("G10").Value = "➖ Neutre"
# Row 11: BRVM - DISTRIBUTION
$ws1.Range("A11").Value = "BRVM - DISTRIBUTION"
$ws1.Range("B11").Value = 0
$ws1.Range("C11").Value = 4
$ws1.Range("D11").Value = 1494.24
$ws1.Range("E11").Value = 368.2
$ws1.Range("F11").Value = "🟡 Observer"
$ws1.Range("G11").Value = "➖ Neutre"
# Row 12: BRVM - TRANSPORT
$ws1.Range("A12").Value = "BRVM - TRANSPORT"
$ws1.Range("B12").Value = 0
$ws1.Range("C12").Value = 4
$ws1.Range("D12").Value = 1397.66
$ws1.Range("E12").Value = 353.7
$ws1.Range("F12").Value = "🟡 Observer"
$ws1.Range("G12").Value = "➖ Neutre"
# Row 13: BRVM - AGRICULTURE
$ws1.Range("A13").Value = "BRVM - AGRICULTURE"
$ws1.Range("B13").Value = 0
$ws1.Range("C13").Value = 4
$ws1.Range("D13").Value = 1320.32
$ws1.Range("E13").Value = 326.32
$ws1.Range("F13").Value = "🟡 Observer"
$ws1.Range("G13").Value = "➖ Neutre"
# Row 14: BRVM - INDUSTRIE
$ws1.Range("A14").Value = "BRVM - INDUSTRIE"
$ws1.Range("B14").Value = 0
$ws1.Range("C14").Value = 4
$ws1.Range("D14").Value = 806.67
$ws1.Range("E14").Value = 199.93
$ws1.Range("F14").Value = "🟡 Observer"
$ws1.Range("G14").Value = "➖ Neutre"
# Row 15: BRVM - CONSOMMATION DE BASE
$ws1.Range("A15").Value = "BRVM - CONSOMMATION DE BASE"
$ws1.Range("B15").Value = 0
$ws1.Range("C15").Value = 4
$ws1.Range("D15").Value = 707.51
$ws1.Range("E15").Value = 174.92
$ws1.Range("F15").Value = "🟡 Observer"
$ws1.Range("G15").Value = "➖ Neutre"
# Row 16: BRVM-PRINCIPAL
$ws1.Range("A16").Value = "BRVM-PRINCIPAL"
$ws1.Range("B16").Value = 0
$ws1.Range("C16").Value = 4
$ws1.Range("D16").Value = 706.5599999999999
$ws1.Range("E16").Value = 175.6
$ws1.Range("F16").Value = "🟡 Observer"
$ws1.Range("G16").Value = "➖ Neutre"
# Row 17: BRVM - INDUSTRIELS
$ws1.Range("A17").Value = "BRVM - INDUSTRIELS"
$ws1.Range("B17").Value = 0
$ws1.Range("C17").Value = 4
$ws1.Range("D17").Value = 532.41
$ws1.Range("E17").Value = 134.55
$ws1.Range("F17").Value = "🟡 Observer"
$ws1.Range("G17").Value = "➖ Neutre"
# Row 18: BRVM-PRESTIGE
$ws1.Range("A18").Value = "BRVM-PRESTIGE"
$ws1.Range("B18").Value = 0
$ws1.Range("C18").Value = 4
$ws1.Range("D18").Value = 523.6900000000001
$ws1.Range("E18").Value = 128.77
$ws1.Range("F18").Value = "🟡 Observer"
$ws1.Range("G18").Value = "➖ Neutre"
# Row 19: BRVM - FINANCES
$ws1.Range("A19").Value = "BRVM - FINANCES"
$ws1.Range("B19").Value = 0
$ws1.Range("C19").Value = 4
$ws1.Range("D19").Value = 486.92
$ws1.Range("E19").Value = 120.31
$ws1.Range("F19").Value = "🟡 Observer"
$ws1.Range("G19").Value = "➖ Neutre"
# Row 20: BRVM - SERVICES FINANCIERS
$ws1.Range("A20").Value = "BRVM - SERVICES FINANCIERS"
$ws1.Range("B20").Value = 0
$ws1.Range("C20").Value = 4
$ws1.Range("D20").Value = 478.54
$ws1.Range("E20").Value = 118.24
$ws1.Range("F20").Value = "🟡 Observer"
$ws1.Range("G20").Value = "➖ Neutre"
# Row 21: BRVM - ENERGIE
$ws1.Range("A21").Value = "BRVM - ENERGIE"
$ws1.Range("B21").Value = 0
$ws1.Range("C21").Value = 4
$ws1.Range("D21").Value = 443.6
$ws1.Range("E21").Value = 110.14
$ws1.Range("F21").Value = "🟡 Observer"
$ws1.Range("G21").Value = "➖ Neutre"
# Row 22: BRVM - CONSOMMATION DISCRETIONNAIRE
$ws1.Range("A22").Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Range("B22").Value = 0
$ws1.Range("C22").Value = 4
$ws1.Range("D22").Value = 425.22
$ws1.Range("E22").Value = 104.67
$ws1.Range("F22").Value = "🟡 Observer"
$ws1.Range("G22").Value = "➖ Neutre"
# Row 23: BRVM - TELECOMMUNICATIONS
$ws1.Range("A23").Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Range("B23").Value = 0
$ws1.Range("C23").Value = 4
$ws1.Range("D23").Value = 379.41
$ws1.Range("E23").Value = 93.84999999999999
$ws1.Range("F23").Value = "🟡 Observer"
$ws1.Range("G23").Value = "➖ Neutre"
# Row 24: BERNABE CI (BNBC)
$ws1.Range("A24").Value = "BERNABE CI (BNBC)"
$ws1.Range("B24").Value = 3
$ws1.Range("C24").Value = 1
$ws1.Range("D24").Value = 8.74
$ws1.Range("E24").Value = -7.14
$ws1.Range("F24").Value = "🟢 Achat"
$ws1.Range("G24").Value = "✅ Renforcer"
# Row 25: FILTISAC CI (FTSC)
$ws1.Range("A25").Value = "FILTISAC CI (FTSC)"
$ws1.Range("B25").Value = 2
$ws1.Range("C25").Value = 0
$ws1.Range("D25").Value = 8.35
$ws1.Range("E25").Value = 0.88
$ws1.Range("F25").Value = "🟡 Observer"
$ws1.Range("G25").Value = "➖ Neutre"
# Row 26: BANK OF AFRICA ML (BOAM)
$ws1.Range("A26").Value = "BANK OF AFRICA ML (BOAM)"
$ws1.Range("B26").Value = 1
$ws1.Range("C26").Value = 0
$ws1.Range("D26").Value = 6.79
$ws1.Range("E26").Value = 6.79
$ws1.Range("F26").Value = "🟡 Observer"
$ws1.Range("G26").Value = "➖ Neutre"
# Row 27: BANK OF AFRICA NG (BOAN)
$ws1.Range("A27").Value = "BANK OF AFRICA NG (BOAN)"
$ws1.Range("B27").Value = 1
$ws1.Range("C27").Value = 1
$ws1.Range("D27").Value = 4.58
$ws1.Range("E27").Value = 5.83
$ws1.Range("F27").Value = "🟡 Observer"
$ws1.Range("G27").Value = "👀 À surveiller"
# Row 28: SOLIBRA CI (SLBC)
$ws1.Range("A28").Value = "SOLIBRA CI (SLBC)"
$ws1.Range("B28").Value = 2
$ws1.Range("C28").Value = 1
$ws1.Range("D28").Value = 3.96
$ws1.Range("E28").Value = -7.48
$ws1.Range("F28").Value = "🟡 Observer"
$ws1.Range("G28").Value = "👀 À surveiller"
# Row 29: NSIA BANQUE COTE D'IVOIRE (NSBC)
$ws1.Range("A29").Value = "NSIA BANQUE COTE D'IVOIRE (NSBC)"
$ws1.Range("B29").Value = 1
$ws1.Range("C29").Value = 0
$ws1.Range("D29").Value = 3.3
$ws1.Range("E29").Value = 3.3
$ws1.Range("F29").Value = "🟡 Observer"
$ws1.Range("G29").Value = "➖ Neutre"
# Row 30: SAFCA CI (SAFC)
$ws1.Range("A30").Value = "SAFCA CI (SAFC)"
$ws1.Range("B30").Value = 1
$ws1.Range("C30").Value = 0
$ws1.Range("D30").Value = 3.08
$ws1.Range("E30").Value = 3.08
$ws1.Range("F30").Value = "🟡 Observer"
$ws1.Range("G30").Value = "➖ Neutre"
# Row 31: SERVAIR ABIDJAN CI (ABJC)
$ws1.Range("A31").Value = "SERVAIR ABIDJAN CI (ABJC)"
$ws1.Range("B31").Value = 1
$ws1.Range("C31").Value = 0
$ws1.Range("D31").Value = 2.99
$ws1.Range("E31").Value = 2.99
$ws1.Range("F31").Value = "🟡 Observer"
$ws1.Range("G31").Value = "➖ Neutre"
# Row 32: SONATEL SN (SNTS)
$ws1.Range("A32").Value = "SONATEL SN (SNTS)"
$ws1.Range("B32").Value = 1
$ws1.Range("C32").Value = 1
$ws1.Range("D32").Value = 1.76
$ws1.Range("E32").Value = -1.83
$ws1.Range("F32").Value = "🟡 Observer"
$ws1.Range("G32").Value = "👀 À surveiller"
# Row 33: CORIS BANK INTERNATIONAL (CBIBF)
$ws1.Range("A33").Value = "CORIS BANK INTERNATIONAL (CBIBF)"
$ws1.Range("B33").Value = 1
$ws1.Range("C33").Value = 1
$ws1.Range("D33").Value = 1.49
$ws1.Range("E33").Value = -5.69
$ws1.Range("F33").Value = "🟡 Observer"
$ws1.Range("G33").Value = "👀 À surveiller"
# Row 34: TOTALENERGIES MARKETING SN (TTLS)
$ws1.Range("A34").Value = "TOTALENERGIES MARKETING SN (TTLS)"
$ws1.Range("B34").Value = 1
$ws1.Range("C34").Value = 0
$ws1.Range("D34").Value = 0.74
$ws1.Range("E34").Value = 0.74
$ws1.Range("F34").Value = "🟡 Observer"
$ws1.Range("G34").Value = "➖ Neutre"
# Row 35: SICABLE CI (CABC)
$ws1.Range("A35").Value = "SICABLE CI (CABC)"
$ws1.Range("B35").Value = 1
$ws1.Range("C35").Value = 1
$ws1.Range("D35").Value = 0.27
$ws1.Range("E35").Value = 5.31
$ws1.Range("F35").Value = "🟡 Observer"
$ws1.Range("G35").Value = "👀 À surveiller"
# Row 36: TRACTAFRIC MOTORS CI (PRSC)
$ws1.Range("A36").Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws1.Range("B36").Value = 1
$ws1.Range("C36").Value = 1
$ws1.Range("D36").Value = 0.15
$ws1.Range("E36").Value = 4
$ws1.Range("F36").Value = "🟡 Observer"
$ws1.Range("G36").Value = "👀 À surveiller"
# Row 37: TOTAL
$ws1.Range("A37").Value = "TOTAL"
$ws1.Range("B37").Value = 0
$ws1.Range("C37").Value = 4
$ws1.Range("D37").Value = 0
$ws1.Range("E37").Value = 0
$ws1.Range("F37").Value = "🟡 Observer"
$ws1.Range("G37").Value = "➖ Neutre"
# Row 38: AFRICA GLOBAL LOGISTICS CI (SDSC)
$ws1.Range("A38").Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$ws1.Range("B38").Value = 1
$ws1.Range("C38").Value = 1
$ws1.Range("D38").Value = 0
$ws1.Range("E38").Value = 1.4
$ws1.Range("F38").Value = "🟡 Observer"
$ws1.Range("G38").Value = "👀 À surveiller"
# Row 39: NEI-CEDA CI (NEIC)
$ws1.Range("A39").Value = "NEI-CEDA CI (NEIC)"
$ws1.Range("B39").Value = 1
$ws1.Range("C39").Value = 1
$ws1.Range("D39").Value = -0.74
$ws1.Range("E39").Value = 2.59
$ws1.Range("F39").Value = "🟡 Observer"
$ws1.Range("G39").Value = "👀 À surveiller"
# Row 40: SOGB CI (SOGC)
$ws1.Range("A40").Value = "SOGB CI (SOGC)"
$ws1.Range("B40").Value = 0
$ws1.Range("C40").Value = 1
$ws1.Range("D40").Value = -1.22
$ws1.Range("E40").Value = -1.22
$ws1.Range("F40").Value = "🟡 Observer"
$ws1.Range("G40").Value = "➖ Neutre"
# Row 41: ECOBANK COTE D''IVOIRE (ECOC)
$ws1.Range("A41").Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$ws1.Range("B41").Value = 1
$ws1.Range("C41").Value = 1
$ws1.Range("D41").Value = -1.82
$ws1.Range("E41").Value = -5.08
$ws1.Range("F41").Value = "🟡 Observer"
$ws1.Range("G41").Value = "👀 À surveiller"
# Row 42: VIVO ENERGY CI (SHEC)
$ws1.Range("A42").Value = "VIVO ENERGY CI (SHEC)"
$ws1.Range("B42").Value = 0
$ws1.Range("C42").Value = 1
$ws1.Range("D42").Value = -1.88
$ws1.Range("E42").Value = -1.88
$ws1.Range("F42").Value = "🟡 Observer"
$ws1.Range("G42").Value = "➖ Neutre"
# Row 43: BANK OF AFRICA BN (BOAB)
$ws1.Range("A43").Value = "BANK OF AFRICA BN (BOAB)"
$ws1.Range("B43").Value = 0
$ws1.Range("C43").Value = 1
$ws1.Range("D43").Value = -2
$ws1.Range("E43").Value = -2
$ws1.Range("F43").Value = "🟡 Observer"
$ws1.Range("G43").Value = "➖ Neutre"
# Row 44: CFAO MOTORS CI (CFAC)
$ws1.Range("A44").Value = "CFAO MOTORS CI (CFAC)"
$ws1.Range("B44").Value = 0
$ws1.Range("C44").Value = 1
$ws1.Range("D44").Value = -2.21
$ws1.Range("E44").Value = -2.21
$ws1.Range("F44").Value = "🟡 Observer"
$ws1.Range("G44").Value = "➖ Neutre"
# Row 45: SMB CI (SMBC)
$ws1.Range("A45").Value = "SMB CI (SMBC)"
$ws1.Range("B45").Value = 0
$ws1.Range("C45").Value = 1
$ws1.Range("D45").Value = -2.48
$ws1.Range("E45").Value = -2.48
$ws1.Range("F45").Value = "🟡 Observer"
$ws1.Range("G45").Value = "➖ Neutre"
# Row 46: BICI CI (BICC)
$ws1.Range("A46").Value = "BICI CI (BICC)"
$ws1.Range("B46").Value = 0
$ws1.Range("C46").Value = 1
$ws1.Range("D46").Value = -2.6
$ws1.Range("E46").Value = -2.6
$ws1.Range("F46").Value = "🟡 Observer"
$ws1.Range("G46").Value = "➖ Neutre"
# Row 47: BANK OF AFRICA BF (BOABF)
$ws1.Range("A47").Value = "BANK OF AFRICA BF (BOABF)"
$ws1.Range("B47").Value = 0
$ws1.Range("C47").Value = 1
$ws1.Range("D47").Value = -2.86
$ws1.Range("E47").Value = -2.86
$ws1.Range("F47").Value = "🟡 Observer"
$ws1.Range("G47").Value = "➖ Neutre"
# Row 48: BANK OF AFRICA CI (BOAC)
$ws1.Range("A48").Value = "BANK OF AFRICA CI (BOAC)"
$ws1.Range("B48").Value = 0
$ws1.Range("C48").Value = 1
$ws1.Range("D48").Value = -3.27
$ws1.Range("E48").Value = -3.27
$ws1.Range("F48").Value = "🟡 Observer"
$ws1.Range("G48").Value = "➖ Neutre"
# Row 49: SETAO CI (STAC)
$ws1.Range("A49").Value = "SETAO CI (STAC)"
$ws1.Range("B49").Value = 0
$ws1.Range("C49").Value = 1
$ws1.Range("D49").Value = -6.09
$ws1.Range("E49").Value = -6.09
$ws1.Range("F49").Value = "🟡 Observer"
$ws1.Range("G49").Value = "➖ Neutre"
# Row 50: UNIWAX CI (UNXC)
$ws1.Range("A50").Value = "UNIWAX CI (UNXC)"
$ws1.Range("B50").Value = 0
$ws1.Range("C50").Value = 1
$ws1.Range("D50").Value = -6.9
$ws1.Range("E50").Value = -6.9
$ws1.Range("F50").Value = "🟡 Observer"
$ws1.Range("G50").Value = "➖ Neutre"

# --- Sheet "Top_YTD": refreshed YTD progression figures (col A unchanged) ---
$ws2.Range("B2").Value = 9573643.18  # BRVM - SERVICES PUBLICS
$ws2.Range("B3").Value = 1373015.24  # SUCRIVOIRE
$ws2.Range("B4").Value = 374504  # SAFCA CI
$ws2.Range("B5").Value = 372225.5  # CFAO MOTORS CI
$ws2.Range("B6").Value = 338759.63  # BRVM - AUTRES SECTEURS
$ws2.Range("B7").Value = 224810  # NEI-CEDA CI
$ws2.Range("B8").Value = 221655.73  # UNIWAX CI
$ws2.Range("B9").Value = 202869.8  # SETAO CI
$ws2.Range("B10").Value = 163760.48  # AIR LIQUIDE CI
$ws2.Range("B11").Value = 50187.55  # BRVM - DISTRIBUTION
